$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Views")

# Mark the XS Landscape (D) / XS Portrait (E) columns as overridden ("Yes")
# for the rows that now have specified overrides. Rows 4 and 17 get a
# special "?Yes" marker in the Portrait column, row 10 only gets the
# Landscape column set.
$ws.Range("D4").Value = "Yes"
$ws.Range("E4").Value = "?Yes"

$ws.Range("D5").Value = "Yes"
$ws.Range("E5").Value = "Yes"

$ws.Range("D6").Value = "Yes"
$ws.Range("E6").Value = "Yes"

$ws.Range("D7").Value = "Yes"
$ws.Range("E7").Value = "Yes"

$ws.Range("D8").Value = "Yes"
$ws.Range("E8").Value = "Yes"

$ws.Range("D9").Value = "Yes"
$ws.Range("E9").Value = "Yes"

$ws.Range("D10").Value = "Yes"

$ws.Range("D11").Value = "Yes"
$ws.Range("E11").Value = "Yes"

$ws.Range("D16").Value = "Yes"
$ws.Range("E16").Value = "Yes"

$ws.Range("D17").Value = "Yes"
$ws.Range("E17").Value = "?Yes"

$ws.Range("D26").Value = "Yes"
$ws.Range("E26").Value = "Yes"

$ws.Range("D29").Value = "Yes"
$ws.Range("E29").Value = "Yes"

# Update the active selection to match the last edited cells
$ws.Activate()
$ws.Range("D29:E29").Select()
